$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("NameTextEntities")
$ws2 = $wb.Worksheets.Item("DescriptionTextEntities")

# Insert a new row before row 10 (the achievement_100 row) on both sheets,
# which pushes that row down to row 11 and inherits the style (s="3",
# "Neutral") of the row above, matching the other "item" rows.
$ws2.Rows.Item(10).Insert()
$ws1.Rows.Item(10).Insert()

# Fill in the new row 10 data for the "Coin" item (key 10304).
# Populate the shared strings in the same order the original author did,
# so new unique strings land at the same indices (41..44) as the target:
#   41 "A coin."   42 "コイン"   43 "Coin"   44 "ただのコイン"
$ws2.Cells.Item(10, 3).Value = "A coin."
$ws1.Cells.Item(10, 4).Value = "コイン"
$ws1.Cells.Item(10, 3).Value = "Coin"
$ws2.Cells.Item(10, 4).Value = "ただのコイン"

# id + key columns on both sheets.
$ws1.Cells.Item(10, 1).Value = 9
$ws1.Cells.Item(10, 2).Value = 10304
$ws2.Cells.Item(10, 1).Value = 9
$ws2.Cells.Item(10, 2).Value = 10304

# The achievement row that got pushed down to row 11 now has id 10.
$ws1.Cells.Item(11, 1).Value = 10
$ws2.Cells.Item(11, 1).Value = 10

# Sheet-view / selection changes: the active tab moves from
# DescriptionTextEntities to NameTextEntities, selections move too.
$ws2.Range("D10").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("D14").Select() | Out-Null

$wb.Save() | Out-Null
